$wb = $excel.ActiveWorkbook

# --- Text change: "Ready for handoff" -> "In Translation" --------------------
# This status string shows up on every sheet that tracks a per-file/per-language
# handoff status: the "Overview" sheet (one status column per language) and the
# per-language sheets ("zh-cn", "de-de") each with their own "Status" column.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Column width changes -----------------------------------------------------
# The status columns got narrower in the regenerated report. Excel stores
# column widths as a pixel-quantized "character width" (MDW-7 grid), so the
# ColumnWidth value below is the character width that lands in the same
# pixel bucket as the archived report's column width.
$newStatusColumnWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColumnWidth  # column E ("zh-cn" status)
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColumnWidth  # column F ("de-de" status)

$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColumnWidth      # column C ("Status")
$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColumnWidth      # column C ("Status")
